# Fix priors: the "Male" block's Linf/k/t0 prior type cells (D18:D20)
# were incorrectly set to "Normal" and should be "No prior", matching
# the "Female" block above (D6:D8) and the CV rows (D9/D21).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D18").Value = "No prior"
$ws.Range("D19").Value = "No prior"
$ws.Range("D20").Value = "No prior"

# Update the active selection to match the saved workbook state.
$ws.Range("D20").Select()
